# Auto-applies the cell-value updates described by the OOXML diff
# (commit: "Update gh-pages to output generated at 456a3b4").
$wb = $excel.ActiveWorkbook

# --- 展览 ---
$ws = $wb.Worksheets.Item('展览')
$ws.Range('F4').Value = 606
$ws.Range('F5').Value = 2659
$ws.Range('F7').Value = 195
$ws.Range('F9').Value = 259
$ws.Range('F10').Value = 5851
$ws.Range('F11').Value = 894
$ws.Range('F13').Value = 1449
$ws.Range('F15').Value = 7012
$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = '2024-08-10'
$ws.Range('C16').Value = '上海·Key Only同人茶会'
$ws.Range('D16').Value = '淞虹路938号福缘湾九华商业广场F1 THOUSE艺术空间'
$ws.Range('E16').Value = '2024.08.10 09:30-08.10 20:00'
$ws.Range('F16').Value = 61
$ws.Range('G16').Value = 120
$ws.Range('H16').Value = 'https://show.bilibili.com/platform/detail.html?id=90039'
$ws.Range('I16').Value = '//i2.hdslb.com/bfs/openplatform/202407/MtZ5eU6D1722321100952.jpeg'
$ws.Range('C17').Value = '上海·创造力动漫游戏嘉年华-风袖，小忻双人内场'
$ws.Range('D17').Value = '莘福路288号 美莘商业广场'
$ws.Range('E17').Value = '2024.08.10 10:00-08.10 17:00'
$ws.Range('F17').Value = 69
$ws.Range('G17').Value = 188
$ws.Range('H17').Value = 'https://show.bilibili.com/platform/detail.html?id=88106'
$ws.Range('I17').Value = '//i2.hdslb.com/bfs/openplatform/202406/IbLKxBuu1719389910566.png'
$ws.Range('C18').Value = '上海·创造力动漫游戏嘉年华1.0'
$ws.Range('E18').Value = '2024.08.10 10:00-08.11 17:00'
$ws.Range('F18').Value = 4830
$ws.Range('G18').Value = 65
$ws.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=87667'
$ws.Range('I18').Value = '//i1.hdslb.com/bfs/openplatform/202406/WRzn64lS1719391076920.png'
$ws.Range('C19').Value = '上海·创造力动漫游戏嘉年华签售票-爱拍照的玉老师'
$ws.Range('E19').Value = '2024.08.10 10:00-08.10 17:00'
$ws.Range('F19').Value = 379
$ws.Range('G19').Value = '已售罄'
$ws.Range('H19').Value = 'https://show.bilibili.com/platform/detail.html?id=88465'
$ws.Range('I19').Value = '//i2.hdslb.com/bfs/openplatform/202407/h6LllgLT1719854666056.png'
$ws.Range('C20').Value = '上海·动漫水着嘉年华'
$ws.Range('D20').Value = '民府路678号 抖音江湾广场'
$ws.Range('E20').Value = '2024.08.10 11:00-08.10 17:00'
$ws.Range('F20').Value = 84
$ws.Range('G20').Value = 68
$ws.Range('H20').Value = 'https://show.bilibili.com/platform/detail.html?id=89929'
$ws.Range('I20').Value = '//i2.hdslb.com/bfs/openplatform/202407/UTF8WBkE1722219440039.jpeg'
$ws.Range('C21').Value = '上海·坏孩纸物语の第48届动漫节之梦回春秋战国（免费活动）（取消）'
$ws.Range('D21').Value = '世纪大道2002号 S.C.Plaza'
$ws.Range('E21').Value = '2024.08.10 10:00-08.11 17:00'
$ws.Range('F21').Value = 743
$ws.Range('G21').Value = '不可售'
$ws.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=88004'
$ws.Range('I21').Value = '//i2.hdslb.com/bfs/openplatform/202406/PP9QKg0v1719293500048.png'
$ws.Range('F22').Value = 2444
$ws.Range('F23').Value = 1302
$ws.Range('F24').Value = 474
$ws.Range('F25').Value = 1180
$ws.Range('F26').Value = 254
$ws.Range('F27').Value = 103
$ws.Range('F28').Value = 107
$ws.Range('F29').Value = 192
$ws.Range('F31').Value = 1320
$ws.Range('F32').Value = 2032
$ws.Range('F33').Value = 265
$ws.Range('F34').Value = 547
$ws.Range('F35').Value = 36
$ws.Range('F37').Value = 1423
$ws.Range('F39').Value = 103
$ws.Range('F40').Value = 539
$ws.Range('F41').Value = 209
$ws.Range('F42').Value = 1678
$ws.Range('F43').Value = 2461
$ws.Range('F45').Value = 94
$ws.Range('F46').Value = 245
$ws.Range('F48').Value = 48

# --- 演出 ---
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F9').Value = 293
$ws.Range('G9').Value = 169
$ws.Range('F11').Value = 70
$ws.Range('F15').Value = 238
$ws.Range('F17').Value = 39
$ws.Range('F24').Value = 311
$ws.Range('F36').Value = 9

# --- 本地生活 ---
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F4').Value = 504
$ws.Range('F6').Value = 1680
$ws.Range('F7').Value = 557
$ws.Range('F8').Value = 1386
$ws.Range('F10').Value = 1784
$ws.Range('F11').Value = 2323
$ws.Range('F12').Value = 750
$ws.Range('F13').Value = 634

# --- 全部类型 ---
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F2').Value = 504
$ws.Range('F4').Value = 606
$ws.Range('F5').Value = 557
$ws.Range('F6').Value = 2659
$ws.Range('F7').Value = 195
$ws.Range('F8').Value = 1386
$ws.Range('F9').Value = 259
$ws.Range('F10').Value = 2323
$ws.Range('F11').Value = 5851
$ws.Range('F12').Value = 750
$ws.Range('F15').Value = 61
$ws.Range('F16').Value = 4831
$ws.Range('F17').Value = 2444
$ws.Range('F18').Value = 1302
$ws.Range('F19').Value = 474
$ws.Range('F20').Value = 1180
$ws.Range('F21').Value = 107
$ws.Range('F22').Value = 293
$ws.Range('G22').Value = 169
$ws.Range('F24').Value = 192
$ws.Range('F25').Value = 70
$ws.Range('F28').Value = 2032
$ws.Range('F29').Value = 265
$ws.Range('F30').Value = 547
$ws.Range('F33').Value = 39
$ws.Range('F34').Value = 1423
$ws.Range('F36').Value = 103
$ws.Range('F37').Value = 539
$ws.Range('F40').Value = 209
$ws.Range('F42').Value = 1679
$ws.Range('F43').Value = 2461
$ws.Range('F44').Value = 94
$ws.Range('F45').Value = 245
$ws.Range('F47').Value = 48
$ws.Range('F49').Value = 9

